# Fruta / hortaliza, semanal
# Weekly refresh of the Membrillo price sheet: rotates the weekly price
# records across rows 2, 3, 5, 8, 9, 10, 11 (adds the newest week's data,
# drops the oldest week forward one slot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44299
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10500
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región del Maule'
$ws.Range("S2").Value = 583
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44299
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región del Maule'
$ws.Range("S3").Value = 500
$ws.Range("T3").Value = 18

# Row 5
$ws.Range("D5").Value = 44425
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("S5").Value = 694

# Row 8
$ws.Range("D8").Value = 44316
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 9500
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 528

# Row 9
$ws.Range("D9").Value = 44272
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 9500
$ws.Range("Q9").Value = '$/caja 15 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 633
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44272
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("S10").Value = 533
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44358
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 639
